$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: test #1 "VERIFICAR RESPONSIVE A FORMATO TABLET" ---
# E4 ("OK") is cleared; F4 stays marked "x" (no textual change, only shared-string
# index shift happens automatically on save).
$ws.Range("E4").ClearContents()

# --- Row 5: test #2 comment changes ---
$ws.Range("G5").Value = "pendiente, ver texto de card y header"

# --- Row 11: test #8 comment changes ---
$ws.Range("G11").Value = "Por ahora hay una sóla Pagina"

# --- Row 13: test #10 "ciudad" — result flips from ERROR/x to OK ---
$ws.Range("E13").Value = "OK"
$ws.Range("F13").ClearContents()
$ws.Range("G13").ClearContents()

# --- Row 32: test #29 "vehiculos" — result flips from ERROR/x to OK ---
$ws.Range("E32").Value = "OK"
$ws.Range("F32").ClearContents()
$ws.Range("G32").ClearContents()

# --- Row 35: test #32 — result flips from ERROR/x to OK ---
$ws.Range("E35").Value = "OK"
$ws.Range("F35").ClearContents()
$ws.Range("G35").ClearContents()

# --- Update the sheet's active selection from F13 to G13 ---
$ws.Activate()
$ws.Range("G13").Select()
